$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C10 ("Integer min" value for the R20 rule row) changes from 18 to 1.
$ws.Range("C10").Value = 1
